# Daily attendance processing - 2025-10-09 18:26:47
# Normalizes the "Recorded By" column (G) so that multi-author entries
# are reordered with the last contributor moved to the front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reordered = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
            $newVal = [string]::Join(", ", $reordered)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
